$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "80.818.59"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "  +5.93%  "

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "3.225.58"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "  +6.12%  "

# Row 4
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "213.38"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "  +7.53%  "

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "636.98"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.97%  "

# Row 7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.273"
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "  +31.03%  "

# Row 8
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.05%  "

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.608"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "  +10.71%  "

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "3.226.68"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "  +6.21%  "

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.629"
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "  +43.53%  "

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000275"
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "  +42.93%  "

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "5.45"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "  +4.94%  "

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "3.817.87"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "  +6.19%  "

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "32.94"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "  +14.06%  "

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "80.661.33"
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "  +5.90%  "

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "3.221.77"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "  +6.44%  "

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "14.66"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "  +8.83%  "

# Row 20
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = "  +25.63%  "

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "9.43"
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "  +5.29%  "

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "448.63"
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "  +17.88%  "

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "5.41"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "  +23.49%  "

# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "4.89"
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "  +12.86%  "

# Row 25
$c = $ws.Cells.Item(25, 2)
$c.NumberFormat = "@"
$c.Value = "WrappedeETH"
$c = $ws.Cells.Item(25, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "3.387.12"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "  +6.58%  "

# Row 26
$c = $ws.Cells.Item(26, 2)
$c.NumberFormat = "@"
$c.Value = "Litecoin"
$c = $ws.Cells.Item(26, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "78.15"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "  +7.87%  "

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "10.99"
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = "  +11.66%  "

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000127"
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = "  +18.28%  "

# Row 29
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.15%  "

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "9.38"
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = "  +13.32%  "

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.11%  "

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "571.14"
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = "  +16.22%  "

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "1.54"
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = "  +10.20%  "

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.160"
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = "  +29.51%  "

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "2.05"
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = "  +6.67%  "

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "23.83"
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = "  +15.71%  "

# Row 37
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = "  +21.05%  "

# Row 38
$c = $ws.Cells.Item(38, 2)
$c.NumberFormat = "@"
$c.Value = "PolygonEcosystemToken"
$c = $ws.Cells.Item(38, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.419"
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = "  +11.37%  "

# Row 39
$c = $ws.Cells.Item(39, 2)
$c.NumberFormat = "@"
$c.Value = "FirstDigitalUSD"
$c = $ws.Cells.Item(39, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "5.88"
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "  +14.88%  "

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "164.41"
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.18%  "

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "20.34"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.46%  "

# Row 43
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "193.12"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.04%  "

# Row 44
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.00%  "

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "1.85"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "  +12.81%  "

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "2.77"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "  +13.94%  "

# Row 47
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "  +9.39%  "

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.808"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.92%  "

# Row 49
$c = $ws.Cells.Item(49, 2)
$c.NumberFormat = "@"
$c.Value = "Filecoin"
$c = $ws.Cells.Item(49, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "4.39"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "  +13.85%  "

# Row 50
$c = $ws.Cells.Item(50, 2)
$c.NumberFormat = "@"
$c.Value = "OKB"
$c = $ws.Cells.Item(50, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "43.64"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "  +6.09%  "

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.650"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "  +9.62%  "
